{"js": "// Replace the three-digit \u00f7 one-digit division problems/answers in the\n// table cells with the new values from the commit.\n// Each original string is unique in the document, so a scoped\n// `body.search(...).insertText(..., \"Replace\")` per pair is safe and\n// exact (no accidental partial/duplicate matches).\n\nconst replacements = [\n  [\"365\u00f78=45, 5\", \"752\u00f78=94, 0\"],\n  [\"496\u00f73=165, 1\", \"698\u00f74=174, 2\"],\n  [\"820\u00f74=205, 0\", \"298\u00f73=99, 1\"],\n  [\"184\u00f74=46, 0\", \"719\u00f76=119, 5\"],\n  [\"232\u00f76=38, 4\", \"928\u00f78=116, 0\"],\n  [\"168\u00f77=24, 0\", \"509\u00f74=127, 1\"],\n  [\"780\u00f79=86, 6\", \"785\u00f76=130, 5\"],\n  [\"653\u00f79=72, 5\", \"356\u00f73=118, 2\"],\n  [\"761\u00f76=126, 5\", \"164\u00f77=23, 3\"],\n  [\"397\u00f73=132, 1\", \"932\u00f77=133, 1\"],\n  [\"761\u00f74=190, 1\", \"282\u00f79=31, 3\"],\n  [\"321\u00f74=80, 1\", \"120\u00f76=20, 0\"],\n  [\"203\u00f76=33, 5\", \"887\u00f74=221, 3\"],\n  [\"104\u00f78=13, 0\", \"125\u00f73=41, 2\"],\n  [\"115\u00f73=38, 1\", \"557\u00f79=61, 8\"],\n  [\"889\u00f78=111, 1\", \"419\u00f74=104, 3\"],\n  [\"949\u00f73=316, 1\", \"750\u00f75=150, 0\"],\n  [\"124\u00f76=20, 4\", \"203\u00f74=50, 3\"],\n  [\"909\u00f73=303, 0\", \"184\u00f73=61, 1\"],\n  [\"234\u00f76=39, 0\", \"176\u00f77=25, 1\"],\n  [\"198\u00f76=33, 0\", \"871\u00f75=174, 1\"],\n  [\"456\u00f77=65, 1\", \"667\u00f75=133, 2\"],\n  [\"352\u00f79=39, 1\", \"536\u00f75=107, 1\"],\n  [\"417\u00f76=69, 3\", \"654\u00f78=81, 6\"],\n  [\"749\u00f74=187, 1\", \"744\u00f76=124, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit \u00f7 one-digit division problems/answers in the\n# table cells with the new values from the commit.\n# Each original string is unique in the document, so a Find/Replace pass\n# per pair (scoped to the whole document content) is exact.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"365\u00f78=45, 5\", \"752\u00f78=94, 0\"),\n    @(\"496\u00f73=165, 1\", \"698\u00f74=174, 2\"),\n    @(\"820\u00f74=205, 0\", \"298\u00f73=99, 1\"),\n    @(\"184\u00f74=46, 0\", \"719\u00f76=119, 5\"),\n    @(\"232\u00f76=38, 4\", \"928\u00f78=116, 0\"),\n    @(\"168\u00f77=24, 0\", \"509\u00f74=127, 1\"),\n    @(\"780\u00f79=86, 6\", \"785\u00f76=130, 5\"),\n    @(\"653\u00f79=72, 5\", \"356\u00f73=118, 2\"),\n    @(\"761\u00f76=126, 5\", \"164\u00f77=23, 3\"),\n    @(\"397\u00f73=132, 1\", \"932\u00f77=133, 1\"),\n    @(\"761\u00f74=190, 1\", \"282\u00f79=31, 3\"),\n    @(\"321\u00f74=80, 1\", \"120\u00f76=20, 0\"),\n    @(\"203\u00f76=33, 5\", \"887\u00f74=221, 3\"),\n    @(\"104\u00f78=13, 0\", \"125\u00f73=41, 2\"),\n    @(\"115\u00f73=38, 1\", \"557\u00f79=61, 8\"),\n    @(\"889\u00f78=111, 1\", \"419\u00f74=104, 3\"),\n    @(\"949\u00f73=316, 1\", \"750\u00f75=150, 0\"),\n    @(\"124\u00f76=20, 4\", \"203\u00f74=50, 3\"),\n    @(\"909\u00f73=303, 0\", \"184\u00f73=61, 1\"),\n    @(\"234\u00f76=39, 0\", \"176\u00f77=25, 1\"),\n    @(\"198\u00f76=33, 0\", \"871\u00f75=174, 1\"),\n    @(\"456\u00f77=65, 1\", \"667\u00f75=133, 2\"),\n    @(\"352\u00f79=39, 1\", \"536\u00f75=107, 1\"),\n    @(\"417\u00f76=69, 3\", \"654\u00f78=81, 6\"),\n    @(\"749\u00f74=187, 1\", \"744\u00f76=124, 0\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
